$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PUBLINDAL")

# Record the "PRUEBA" description for the existing accident report (row 2)
$ws.Range("G2").Value = "PRUEBA"

# Add the new accident report row (row 3)
$ws.Range("A3").Value = "Accidente con baja"
$ws.Range("B3").Value = "Publindal"

# Force the date to be stored as text so it matches the "2024-11-18" string
# used elsewhere in the sheet instead of being auto-converted to a date serial.
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2024-11-18"

$ws.Range("D3").Value = "14:30"
$ws.Range("E3").Value = "MARÍA GONZALEZ ALARCÓN"
$ws.Range("F3").Value = "SERIGRAFÍA"
$ws.Range("G3").Value = "PRUEBA"
$ws.Range("H3").Value = "MÁQUINA SERIGRAFÍA"
$ws.Range("I3").Value = "TORAX"
$ws.Range("J3").Value = "CORTE"
